# Rename Sheet1 to WC_2022
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "WC_2022"

# Update the selection to D26 (single cell) to match the saved view state
$ws.Range("D26").Select()
